# ---------------------------------------------------------------------------
# Add a new "2022-Q3" quarter: insert its summary row into "总计" and create
# a brand-new "2022-Q3" worksheet (positioned right after "总计") holding the
# per-fund holdings table for that quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row right under the header for the
#    2022-Q3 totals; existing quarters shift down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 22
$total.Cells.Item(2, 4).Value = 4.58

# Column A carries the bold/bordered/centered style used throughout the
# sheet (copy it from the row directly below, which already has it).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Cells.Item(2, 1).Value = 0

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: fund-level holdings table, placed right after
#    "总计" (i.e. before the existing "2022-Q2" tab).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Match the page margins used throughout the rest of the workbook (Excel's
# PageSetup is in points; 72pt == 1in).
$q3.PageSetup.LeftMargin = 0.75 * 72
$q3.PageSetup.RightMargin = 0.75 * 72
$q3.PageSetup.TopMargin = 1 * 72
$q3.PageSetup.BottomMargin = 1 * 72
$q3.PageSetup.HeaderMargin = 0.5 * 72
$q3.PageSetup.FooterMargin = 0.5 * 72

# -- header row --------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# -- fund rows -----------------------------------------------------------
# index, code, name, size, stockPosition, positionRatio, marketValue, rank
$rows = @(
    @(0,  "007130", "中庚小盘价值股票",           "75.87", "93.06", "3.11", "2.3596", 7),
    @(1,  "007497", "中庚价值灵动灵活配置混合",    "36.46", "89.30", "4.66", "1.6990", 1),
    @(2,  "000780", "鹏华医疗保健股票",           "6.60",  "81.50", "3.10", "0.2046", 10),
    @(3,  "009913", "中信保诚成长动力混合A",       "5.03",  "74.01", "3.11", "0.1564", 7),
    @(4,  "014220", "恒越医疗健康精选混合A",       "0.72",  "88.76", "3.84", "0.0276", 9),
    @(5,  "003284", "中邮医药健康灵活配置混合",    "0.65",  "76.19", "3.22", "0.0209", 6),
    @(6,  "003513", "中邮消费升级灵活配置混合",    "0.56",  "30.56", "3.16", "0.0177", 3),
    @(7,  "001415", "信诚新锐回报灵活配置混合A",   "2.37",  "26.42", "0.61", "0.0145", 9),
    @(8,  "002046", "信诚新锐回报灵活配置混合B",   "2.07",  "26.42", "0.61", "0.0126", 9),
    @(9,  "014221", "恒越医疗健康精选混合C",       "0.29",  "88.76", "3.84", "0.0111", 9),
    @(10, "003235", "信诚至利灵活配置混合C",       "1.82",  "25.21", "0.61", "0.0111", 8),
    @(11, "001402", "信诚新选回报灵活配置混合A",   "1.37",  "23.94", "0.61", "0.0084", 7),
    @(12, "008037", "兴银先锋成长混合A",           "0.21",  "71.76", "3.38", "0.0071", 1),
    @(13, "014285", "鑫元健康产业混合A",           "0.12",  "78.73", "3.75", "0.0045", 1),
    @(14, "008038", "兴银先锋成长混合C",           "0.13",  "71.76", "3.38", "0.0044", 1),
    @(15, "004157", "信诚至诚灵活配置混合A",       "0.63",  "24.85", "0.63", "0.0040", 8),
    @(16, "003234", "信诚至利灵活配置混合A",       "0.66",  "25.21", "0.61", "0.0040", 8),
    @(17, "014286", "鑫元健康产业混合C",           "0.07",  "78.73", "3.75", "0.0026", 1),
    @(18, "002030", "信诚新选回报灵活配置混合B",   "0.38",  "23.94", "0.61", "0.0023", 7),
    @(19, "014282", "中信保诚成长动力混合C",       "0.07",  "74.01", "3.11", "0.0022", 7),
    @(20, "001474", "兴银丰盈灵活配置混合",        "0.09",  "60.49", "1.71", "0.0015", 8),
    @(21, "004158", "信诚至诚灵活配置混合B",       "0.17",  "24.85", "0.63", "0.0011", 8)
)

# Text-valued columns (B..G) must stay text even though they look numeric
# (e.g. fund codes with leading zeros, percentages stored as plain text).
# Force text storage via a temporary "@" format, then ClearFormats() right
# after writing the values so no lingering number-format style survives on
# the cell (matches the source data, which carries no explicit style here).
$q3.Range("B2:G23").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q3.Cells.Item($excelRow, 1).Value = [int]$row[0]
    $q3.Cells.Item($excelRow, 2).Value = [string]$row[1]
    $q3.Cells.Item($excelRow, 3).Value = [string]$row[2]
    $q3.Cells.Item($excelRow, 4).Value = [string]$row[3]
    $q3.Cells.Item($excelRow, 5).Value = [string]$row[4]
    $q3.Cells.Item($excelRow, 6).Value = [string]$row[5]
    $q3.Cells.Item($excelRow, 7).Value = [string]$row[6]
    $q3.Cells.Item($excelRow, 8).Value = [int]$row[7]
}

$q3.Range("B2:G23").ClearFormats()

# -- formatting: bold/bordered/centered header row + index column --------
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$total.Range("A3").Copy()
$q3.Range("A2:A23").PasteSpecial(-4122)

# Restore the originally active tab ("总计") and clear the Office clipboard
# marching-ants selection left over from the Copy() calls above.
$excel.CutCopyMode = 0
$total.Activate()
[void]$total.Range("A1").Select()

